# Apply fixes to PolishingOrders schedule: correct column D (price) break down
# values and update the active sheet view/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new value for column D (the "price" break down fix)
$updates = @{
    2  = 350
    3  = 350
    4  = 350
    5  = 350
    6  = 350
    7  = 350
    8  = 350
    9  = 350
    10 = 350
    11 = 350
    12 = 350
    13 = 350
    14 = 3500
    16 = 3500
    17 = 350
    18 = 350
    19 = 350
    20 = 350
    21 = 350
    22 = 350
    23 = 350
    24 = 350
    25 = 350
    26 = 350
    27 = 350
    28 = 350
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 4).Value = $updates[$row]
}

# Update the sheet view: scroll position and current selection
$window = $excel.ActiveWindow
$window.ScrollRow = 4
$ws.Range("H23").Select()
